$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.427.26'
$ws.Range("E2").Value = '  +7.68%  '
$ws.Range("D3").Value = '3.654.31'
$ws.Range("E3").Value = '  +18.89%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.54%  '
$ws.Range("D7").Value = '3.651.73'
$ws.Range("E7").Value = '  +18.90%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +5.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.56'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.72%  '
$ws.Range("E12").Value = '  +7.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000260'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.85%  '
$ws.Range("D15").Value = '4.266.75'
$ws.Range("E15").Value = '  +19.01%  '
$ws.Range("D16").Value = '71.378.39'
$ws.Range("E16").Value = '  +7.76%  '
$ws.Range("D17").Value = '3.659.39'
$ws.Range("E17").Value = '  +19.20%  '
$ws.Range("E18").Value = '  +2.66%  '
$ws.Range("E19").Value = '  +9.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '515.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +22.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.750'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.49%  '
$ws.Range("E26").Value = '  +9.97%  '
$ws.Range("E27").Value = '  +8.16%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("E29").Value = '  +14.01%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000113'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +23.72%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.29'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +16.92%  '
$ws.Range("E33").Value = '  +7.02%  '
$ws.Range("E34").Value = '  +5.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.92%  '
$ws.Range("E37").Value = '  +8.83%  '
$ws.Range("E38").Value = '  +12.43%  '
$ws.Range("E39").Value = '  +9.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '47.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.23%  '
$ws.Range("E42").Value = '  +5.68%  '
$ws.Range("E43").Value = '  +8.64%  '
$ws.Range("D44").Value = '3.161.28'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '406.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.41%  '
$ws.Range("E47").Value = '  +7.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +16.57%  '
$ws.Range("E49").Value = '  +16.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("E51").Value = '  +0.02%  '
